# Stroop task slide (slide 22): rework into a "LIO" based button/number layout.
#
# Behaviour being reproduced (per the target diff):
#  - The 4 coloured ovals (answer buttons) are each replaced by a freshly
#    created shape at the same position/size/fill, but carrying a new
#    (permuted) digit label:
#       old "Oval 1" (red,    text "1") -> new "Oval 2" (red,    text "3")
#       old "Oval 3" (green,  text "2") -> new "Oval 7" (green,  text "4")
#       old "Oval 4" (none,   text "3") -> new "Oval 8" (none,   text "1")
#       old "Oval 5" (yellow, text "4") -> new "Oval 9" (yellow, text "2")
#  - The instructions textbox ("TextBox 6") is kept (same shape/id) but is
#    moved to the back of the z-order (it ends up first in the shape tree).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)

# Grab references to the 4 existing ovals and the textbox by their current
# (pre-edit) names, before anything else changes the collection around.
$oval1 = $s.Shapes.Item("Oval 1")
$oval3 = $s.Shapes.Item("Oval 3")
$oval4 = $s.Shapes.Item("Oval 4")
$oval5 = $s.Shapes.Item("Oval 5")
$textBox = $s.Shapes.Item("TextBox 6")

# Duplicate each oval (in original order) to create its replacement shape,
# then put the duplicate back exactly where the original was (Duplicate()
# nudges the copy by a small offset) and rename/retext it.
$newOval2 = $oval1.Duplicate()
$newOval2.Left = $oval1.Left
$newOval2.Top = $oval1.Top
$newOval2.Name = "Oval 2"
$newOval2.TextFrame.TextRange.Text = "3"

$newOval7 = $oval3.Duplicate()
$newOval7.Left = $oval3.Left
$newOval7.Top = $oval3.Top
$newOval7.Name = "Oval 7"
$newOval7.TextFrame.TextRange.Text = "4"

$newOval8 = $oval4.Duplicate()
$newOval8.Left = $oval4.Left
$newOval8.Top = $oval4.Top
$newOval8.Name = "Oval 8"
$newOval8.TextFrame.TextRange.Text = "1"

$newOval9 = $oval5.Duplicate()
$newOval9.Left = $oval5.Left
$newOval9.Top = $oval5.Top
$newOval9.Name = "Oval 9"
$newOval9.TextFrame.TextRange.Text = "2"

# Remove the old ovals now that their replacements are in place.
$oval1.Delete()
$oval3.Delete()
$oval4.Delete()
$oval5.Delete()

# Send the instructions textbox to the back of the z-order.
$textBox.ZOrder(1)
